# Separate combined config into data generation and training configs:
# rename config-1.yml -> data-gen-config-001.yml
# rename config-2.yml -> data-gen-config-002.yml
# across the "Dataset Registry", "Configuration Details" and "File Metadata" sheets.

$wb = $excel.ActiveWorkbook

# Sheet: Dataset Registry (Config File column B)
$wsRegistry = $wb.Worksheets.Item("Dataset Registry")
$wsRegistry.Range("B2").Value = "data-gen-config-001.yml"
$wsRegistry.Range("B3").Value = "data-gen-config-002.yml"

# Sheet: Configuration Details (Config File column B)
$wsConfig = $wb.Worksheets.Item("Configuration Details")
$wsConfig.Range("B2").Value = "data-gen-config-001.yml"
$wsConfig.Range("B3").Value = "data-gen-config-002.yml"

# Sheet: File Metadata (Config Path column B, includes directory prefix)
$wsMeta = $wb.Worksheets.Item("File Metadata")
$wsMeta.Range("B2").Value = "configs\data_generation\data-gen-config-001.yml"
$wsMeta.Range("B3").Value = "configs\data_generation\data-gen-config-002.yml"
